$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("analysis")

# --- Update day 15 (row 16) measurements after the performance improvement ---
$ws.Range("B16").Value = [double]"0.22489999999999999"
$ws.Range("C16").Value = [double]"4.0258000000000002E-2"
$ws.Range("D16").Value = [double]"5.4749999999999998E-3"
$ws.Range("E16").Value = [double]"3.6999999999999999E-4"
$ws.Range("F16").Value = [double]"890.40279999999996"
$ws.Range("G16").Value = [double]"25.593371000000001"

# --- Remove the scratch sheets that held intermediate calculations ---
$null = $wb.Worksheets.Item("Sheet4").Delete()
$null = $wb.Worksheets.Item("Sheet3").Delete()
$null = $wb.Worksheets.Item("Sheet2").Delete()
$null = $wb.Worksheets.Item("Sheet1").Delete()

# --- Bold the day index column and the header row ---
$ws.Range("A1:A29").Font.Bold = $true
$ws.Range("B1:G1").Font.Bold = $true

# --- Resize the chart so it spans further to the right ---
$co = $ws.ChartObjects().Item(1)
$co.Width = 1146.6251181102361

# --- Move the active selection ---
$null = $ws.Range("A3").Select()
